$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: split the "Answer:" run under Question 2 into 5 runs:
#   " " | "I believe that " | "the code is" |
#   " efficient, but definitely it can better" |
#   ", to avoid multiple language issues, I first converted the string
#    to Unicode values, then used the rage of emojis Unicode to
#    exclude them."
# ------------------------------------------------------------------

$oldSentence = " the code is not the best, to avoid multiple language issues, I first converted the string to Unicode values, then used the rage of emojis Unicode to exclude them."

$piece1 = " "
$piece2 = "I believe that "
$piece3 = "the code is"
$piece4 = " efficient, but definitely it can better"
$piece5 = ", to avoid multiple language issues, I first converted the string to Unicode values, then used the rage of emojis Unicode to exclude them."

# Step 1: replace the whole sentence with 4-char tokens (keeps this a
# single run with the original run's formatting intact, since
# Find/Replace preserves rPr of the matched text).
$findRange = $d.Content
$findRange.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "TOK1TOK2TOK3TOK4TOK5", 2)

# Step 2: locate the token block and compute the boundaries between
# the five 4-character tokens.
$tokRange = $d.Content
$tokRange.Find.Execute("TOK1TOK2TOK3TOK4TOK5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s = $tokRange.Start

$b1 = $s + 4
$b2 = $s + 8
$b3 = $s + 12
$b4 = $s + 16

# Step 3: drop bookmarks on each boundary so that the runs created in
# step 4 do not get silently re-coalesced with their neighbour just
# because they end up with identical run formatting.
$d.Bookmarks.Add("zzSplit1", $d.Range($b1, $b1)) | Out-Null
$d.Bookmarks.Add("zzSplit2", $d.Range($b2, $b2)) | Out-Null
$d.Bookmarks.Add("zzSplit3", $d.Range($b3, $b3)) | Out-Null
$d.Bookmarks.Add("zzSplit4", $d.Range($b4, $b4)) | Out-Null

# Step 4: replace each token with its final text, one at a time; each
# Find/Replace preserves the run formatting of the token it replaces.
$r1 = $d.Content
$r1.Find.Execute("TOK1", $true, $false, $false, $false, $false, $true, 1, $false, $piece1, 2)
$r2 = $d.Content
$r2.Find.Execute("TOK2", $true, $false, $false, $false, $false, $true, 1, $false, $piece2, 2)
$r3 = $d.Content
$r3.Find.Execute("TOK3", $true, $false, $false, $false, $false, $true, 1, $false, $piece3, 2)
$r4 = $d.Content
$r4.Find.Execute("TOK4", $true, $false, $false, $false, $false, $true, 1, $false, $piece4, 2)
$r5 = $d.Content
$r5.Find.Execute("TOK5", $true, $false, $false, $false, $false, $true, 1, $false, $piece5, 2)

# Step 5: clean up the temporary bookmarks; the runs remain separate.
$d.Bookmarks("zzSplit1").Delete()
$d.Bookmarks("zzSplit2").Delete()
$d.Bookmarks("zzSplit3").Delete()
$d.Bookmarks("zzSplit4").Delete()

# ------------------------------------------------------------------
# Change 2: merge the " " run and the "The main limitation..." run
# (under Question 3) back into a single run.
# ------------------------------------------------------------------

$mergeRange = $d.Content
$mergeRange.Find.Execute(" The main limitation is that the ranges included might not cover all emojis there.", $true, $false, $false, $false, $false, $true, 1, $false, " The main limitation is that the ranges included might not cover all emojis there.", 2)

Write-Output $d.Content.Text
